# Increase axis text size on both charts in the deck.
#
# Slide 1 -> Chart 1 ("FY16".."FY19" line chart)
# Slide 2 -> Chart 2 ("Actual"/"Forecast"/"Upper"/"Lower" line chart)
#
# For each chart:
#   - chart title font size: 18.62 -> 24 pt
#   - category (date) axis tick-label font size: 11.97 -> 14 pt
#   - value axis tick-label font size: 11.97 -> 14 pt
#   - value axis title font size: 13.3 -> 16 pt
#
# Chart 1 additionally had its last series ("FY16") data-label font size
# nudged from 11.97 -> 12 pt as part of the same edit.

$p = $ppt.ActivePresentation

# --- Slide 1 / Chart 1 ---
$slide1 = $p.Slides.Item(1)
$chart1 = $slide1.Shapes.Item(1).Chart

$chart1.ChartTitle.Format.TextFrame2.TextRange.Font.Size = 24

$chart1CatAxis = $chart1.Axes(1)
$chart1CatAxis.TickLabels.Font.Size = 14

$chart1ValAxis = $chart1.Axes(2)
$chart1ValAxis.TickLabels.Font.Size = 14
$chart1ValAxis.AxisTitle.Format.TextFrame2.TextRange.Font.Size = 16

$chart1Series = $chart1.FullSeriesCollection()
$chart1LastSeries = $chart1Series.Item($chart1Series.Count)
$chart1LastSeries.DataLabels().Format.TextFrame2.TextRange.Font.Size = 12

# --- Slide 2 / Chart 2 ---
$slide2 = $p.Slides.Item(2)
$chart2 = $slide2.Shapes.Item(1).Chart

$chart2.ChartTitle.Format.TextFrame2.TextRange.Font.Size = 24

$chart2CatAxis = $chart2.Axes(1)
$chart2CatAxis.TickLabels.Font.Size = 14

$chart2ValAxis = $chart2.Axes(2)
$chart2ValAxis.TickLabels.Font.Size = 14
$chart2ValAxis.AxisTitle.Format.TextFrame2.TextRange.Font.Size = 16
